# Commit: "remove Gamelogic project, modify SLG building config"
#
# Populate column B (the effect-icon / setting-sheet reference column) for
# every existing row on Sheet1. Row 1 already carries the "Atlas_ResID"
# header in B1; row 2 (EFT_INFO) points at "msg_icon", and every remaining
# row (EFT_BOOOST .. EFT_FINISH) points at the shared "Ssetting" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "msg_icon"
$ws.Range("B3:B15").Value = "Ssetting"

# Leave the cursor where the author left it when the workbook was saved.
$ws.Range("E14").Select() | Out-Null
